$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The slide currently has two top-level shapes: the "Rectangle 19" shape
# and the "Group 14" group. The edit groups them together into a new
# enclosing group (equivalent to selecting both and pressing Ctrl+G).
#
# PowerPoint assigns new shape/group IDs as the lowest unused id on the
# slide. In the target file the new group ends up with id 21 (and the
# default resulting name "Group 20"), which is the 10th unused id after
# the ones already present on the slide (1,4,5,7,8,9,10,11,14,15,20).
# To land on that same id (and therefore the same default name) via the
# COM id-assignment algorithm, we briefly create and discard 9 throwaway
# textboxes so the following Group() call consumes id 21.
for ($i = 0; $i -lt 9; $i++) {
    $dummy = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
    $dummy.Delete()
}

$range = $s.Shapes.Range(@(1, 2))
$grouped = $range.Group()
